# New crime data collected - weekly CompStat 104th Precinct update
# Updates: Police Commissioner name, volume/report-week header text,
# and the weekly/28-day/YTD/2-year crime-count figures for rows 14-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates
# ---------------------------------------------------------------------

# Police Commissioner name (M6)
$ws.Range("M6").Value = "Edward A. Caban"

# Volume "Number" within the rich "Volume 30   Number  26" string (A8):
# only the trailing "26" run changes to "27".
$ws.Range("A8").Characters(21, 2).Text = "27"

# Report covering week dates within "Report Covering the Week  6/26/2023
# Through  7/2/2023" (C9): "6/26/2023" -> "7/3/2023" and "7/2/2023" -> "7/9/2023"
$ws.Range("C9").Characters(27, 9).Text = "7/3/2023"
$ws.Range("C9").Characters(46, 8).Text = "7/9/2023"

# ---------------------------------------------------------------------
# Row 14 (Murder): G14/H14 flip from numeric (1 / -100) to the "no data"
# text markers "0" / "***.*" already used by C14:F14 (style 14).
# ---------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))

# ---------------------------------------------------------------------
# Row 15 (Rape)
# ---------------------------------------------------------------------
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = 0
$ws.Range("N15").Value = -31.25

# ---------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = -43.75
$ws.Range("I16").Value = 113
$ws.Range("J16").Value = 106
$ws.Range("K16").Value = 6.603773584905
$ws.Range("L16").Value = 98.245614035087
$ws.Range("M16").Value = -9.6
$ws.Range("N16").Value = -76.985743380855

# ---------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 6.25
$ws.Range("I17").Value = 154
$ws.Range("J17").Value = 152
$ws.Range("K17").Value = 1.315789473684
$ws.Range("L17").Value = 24.193548387096
$ws.Range("M17").Value = 25.203252032520
$ws.Range("N17").Value = 6.944444444444

# ---------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 8
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = -40.625
$ws.Range("I18").Value = 95
$ws.Range("J18").Value = 141
$ws.Range("K18").Value = -32.624113475177
$ws.Range("L18").Value = -22.131147540983
$ws.Range("M18").Value = -59.574468085106
$ws.Range("N18").Value = -91.079812206572

# ---------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 40
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = 5.555555555555
$ws.Range("I19").Value = 344
$ws.Range("J19").Value = 328
$ws.Range("K19").Value = 4.878048780487
$ws.Range("L19").Value = 42.738589211618
$ws.Range("M19").Value = 60.747663551401
$ws.Range("N19").Value = 8.517350157728

# ---------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = -3.703703703703
$ws.Range("I20").Value = 178
$ws.Range("J20").Value = 135
$ws.Range("K20").Value = 31.851851851851
$ws.Range("L20").Value = 64.814814814814
$ws.Range("M20").Value = -7.772020725388
$ws.Range("N20").Value = -90.476190476190

# ---------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -5.128205128205
$ws.Range("F21").Value = 156
$ws.Range("H21").Value = -13.333333333333
$ws.Range("I21").Value = 897
$ws.Range("J21").Value = 875
$ws.Range("K21").Value = 2.514285714285
$ws.Range("L21").Value = 35.294117647058
$ws.Range("M21").Value = -0.111358574610
$ws.Range("N21").Value = -77.070552147239

# ---------------------------------------------------------------------
# Row 22 (Transit): C22 flips from numeric 1 to the text "0" marker
# (matching D22's style/value), F22 drops from 3 to 2.
# ---------------------------------------------------------------------
$ws.Range("D22").Copy($ws.Range("C22"))
$ws.Range("F22").Value = 2

# ---------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -18.518518518518
$ws.Range("G24").Value = 120
$ws.Range("H24").Value = 3.333333333333
$ws.Range("I24").Value = 659
$ws.Range("J24").Value = 729
$ws.Range("K24").Value = -9.602194787379
$ws.Range("L24").Value = -4.768786127167
$ws.Range("M24").Value = 16.021126760563

# ---------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = -30
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 63
$ws.Range("H25").Value = -33.333333333333
$ws.Range("I25").Value = 267
$ws.Range("J25").Value = 281
$ws.Range("K25").Value = -4.982206405693
$ws.Range("L25").Value = 16.593886462882
$ws.Range("M25").Value = -31.713554987212

# ---------------------------------------------------------------------
# Row 26 (UCR Rape*)
# ---------------------------------------------------------------------
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -11.111111111111

# ---------------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = -62.5
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 37
$ws.Range("J27").Value = 41
$ws.Range("K27").Value = -9.756097560975
$ws.Range("L27").Value = 19.354838709677

# ---------------------------------------------------------------------
# Row 28 (Shooting Vic.)
# ---------------------------------------------------------------------
$ws.Range("G28").Value = 1

# ---------------------------------------------------------------------
# Row 29 (Shooting Inc.)
# ---------------------------------------------------------------------
$ws.Range("G29").Value = 1
